$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host ($ws.Cells.Item(1,1).Value2)
Write-Host ($ws.Cells.Item(1,1).Text)
